$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 73 ("189-Rotate Array"),
# pushing all the rows below it down by two. Excel copies the formatting
# of the row above the insertion point (row 72) onto the new rows, which
# matches the wrap-text / style pattern used throughout this table.
$ws.Rows("73:74").Insert()

# --- Row 74: 155-Min Stack (filled in first, matching the order new
# shared strings were authored in) ---
$ws.Range("A74").Value = "Array"
$ws.Range("B74").Value = 155
$ws.Range("C74").Value = "155-Min Stack"
$ws.Range("D74").Value = "Medium"
$ws.Range("E74").Value = "Two stacks, one is the actual stack, the other used to track min values"
$ws.Range("F74").Value = "O(1) time, O(n) memory"
$ws.Range("G74").Value = "O(1) time, O(n) memory"
$ws.Range("H74").Value = "Two stacks, one is the actual stack, the other used to track min values"
$ws.Range("I74").Value = "O(n)"
$ws.Range("J74").Value = "no"
$ws.Range("K74").Value = "no"
$ws.Range("L74").Value = "When I submitted my code, I missed one edge case where the same element could be pushed multiple times. I had to make sure the new minimums were also added to the second stack"
$ws.Range("M74").Value = "40 minutes"
$ws.Rows("74").RowHeight = 130.5

# --- Row 73: 150-Evaluate Reverse Polish Notation ---
$ws.Range("A73").Value = "Stack"
$ws.Range("B73").Value = 150
$ws.Range("C73").Value = "150-Evaluate Reverse Polish Notation"
$ws.Range("D73").Value = "Medium"
$ws.Range("E73").Value = "One pass scan using a stack"
$ws.Range("F73").Value = "O(n) time, O(n) memory"
$ws.Range("G73").Value = "O(n) time, O(n) memory"
$ws.Range("H73").Value = "One pass scan using a stack"
$ws.Range("I73").Value = "O(n) time, O(n) memory"
$ws.Range("J73").Value = "no"
$ws.Range("K73").Value = "no"
$ws.Range("M73").Value = "20 minutes"

# --- Row 79 (the first of the blank spacer rows after the shift) now
# gets filled in with 739-Daily Temperature ---
$ws.Range("A79").Value = "Stack"
$ws.Range("B79").Value = 739
$ws.Range("C79").Value = "739-Daily Temperature"
$ws.Range("D79").Value = "Medium"
$ws.Range("E79").Value = "Linear scan utilizing stack"
$ws.Range("F79").Value = "O(n) time, O(n) memory"
$ws.Range("G79").Value = "O(n) time, O(1) memory"
$ws.Range("H79").Value = "Linear scan utilizing stack with the array in place"
$ws.Range("J79").Value = "no"
$ws.Range("K79").Value = "no"
$ws.Range("M79").Value = "30 minutes"

# Two more blank spacer rows are appended at the bottom (86 & 87), matching
# the style ("Other" column, text number format + wrap text) used by the
# other spacer rows above.
$ws.Range("L86").WrapText = $true
$ws.Range("L86").NumberFormat = "@"
$ws.Range("L87").WrapText = $true
$ws.Range("L87").NumberFormat = "@"

# Resize (expand) the worksheet table to include the 2 new rows in the
# middle plus the 2 new blank rows appended at the end.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A2:X87"))

# Match the author's final selection/active cell.
$ws.Range("E88").Select()
